# Update the two time-slot values in column C (rows 2-3) that were
# re-labelled from the 9:30/9:35 slot to the 2:55/3:00 slot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "2:55-3:0"
$ws.Range("C3").Value = "3:0-3:5"

# Reflect the author's final on-sheet selection (cell B11 was clicked,
# replacing the earlier C10:C11 selection).
[void]$ws.Range("B11").Select()
